# Auto-committed on 2022/02/17 週四
# Updates the FT_L5060 test-case sheet:
#  - bump the referenced spec doc version from V1.64 to V1.65 (column M, all data rows)
#  - relabel E38 from an "output field" to a "button" ([輸出欄位]是否指定 -> [按鈕]是否指定)
#  - change G38's description to point at the new L5607 assignment-maintenance link,
#    and pick up the "link" cell style (same look as the other 連結至【...】 cells)
#  - bump the test-case open date (column Q) for every data row
#  - move the sheet's active-cell selection down to A42

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column M: requirement-spec version reference, rows 2-40 ---
$ws.Range("M2:M40").Value = "製作依據之需求規格書與版本：PJ201800012_URS_5管理性作業_V1.65.DOCX"

# --- Column Q: test-case open date, rows 2-40 (2022/02/10 -> 2022/02/16) ---
$ws.Range("Q2:Q40").Value = 44608

# --- Row 38: relabel E38 and replace G38's text/style ---
$ws.Range("E38").Value = "[按鈕]是否指定"

# Pick up the "link" look (style used by the other 連結至【...】 cells, e.g. G37)
# before overwriting the text, so G38 matches the new content's formatting.
$ws.Range("G37").Copy() | Out-Null
$ws.Range("G38").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("G38").Value = "連結至【L5607個案人員指派維護】"

# --- Move the active selection to A42 ---
$ws.Range("A42").Select() | Out-Null
